$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5500
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -9248

$ws.Range("H64").Value = 3073.8572
$ws.Range("J64").Value = 3120.6924
$ws.Range("L64").Value = 3120.6924
$ws.Range("N64").Value = -3616.6924

$ws.Range("H65").Value = 5500
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -46240

$ws.Range("H67").Value = 3073.8572
$ws.Range("J67").Value = 3120.6924
$ws.Range("L67").Value = 3120.6924
$ws.Range("N67").Value = -4836.6924

$ws.Range("H138").Value = 3093.8765
$ws.Range("I138").Value = 1697.3914
$ws.Range("J138").Value = 3647.6553
$ws.Range("K138").Value = 5092.174199999999
$ws.Range("L138").Value = 10942.9659
$ws.Range("M138").Value = 47.82580000000053
$ws.Range("N138").Value = -21222.9659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 19469.385
$ws.Range("J101").Value = 19469.385
$ws.Range("L101").Value = 19469.385
$ws.Range("N101").Value = -25959.385

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12655.27
$ws.Range("I94").Value = 8753.462
$ws.Range("J94").Value = 16557.076
$ws.Range("K94").Value = 8753.462
$ws.Range("L94").Value = 16557.076
$ws.Range("M94").Value = -8302.462
$ws.Range("N94").Value = -17459.076

$ws.Range("H107").Value = 4762.2144
$ws.Range("I107").Value = 7024.2354
$ws.Range("K107").Value = 7024.2354
$ws.Range("M107").Value = -5104.2354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1601.3485
$ws.Range("I31").Value = 1091.3695
$ws.Range("J31").Value = 2774.3
$ws.Range("K31").Value = 1091.3695
$ws.Range("L31").Value = 2774.3
$ws.Range("M31").Value = -796.3695
$ws.Range("N31").Value = -3364.3

$ws.Range("H34").Value = 1601.3485
$ws.Range("I34").Value = 1091.3695
$ws.Range("J34").Value = 2774.3
$ws.Range("K34").Value = 1091.3695
$ws.Range("L34").Value = 2774.3
$ws.Range("M34").Value = -889.3695
$ws.Range("N34").Value = -3178.3

$ws.Range("H62").Value = 3109.7334
$ws.Range("I62").Value = 2520
$ws.Range("J62").Value = 3404.6
$ws.Range("K62").Value = 2520
$ws.Range("L62").Value = 3404.6
$ws.Range("M62").Value = -1896
$ws.Range("N62").Value = -4652.6

$ws.Range("H65").Value = 3109.7334
$ws.Range("I65").Value = 2520
$ws.Range("J65").Value = 3404.6
$ws.Range("K65").Value = 12600
$ws.Range("L65").Value = 17023
$ws.Range("M65").Value = -9480
$ws.Range("N65").Value = -23263

$ws.Range("H134").Value = 2498.4614
$ws.Range("I134").Value = 2395.7307
$ws.Range("J134").Value = 2909.3845
$ws.Range("K134").Value = 7187.1921
$ws.Range("L134").Value = 8728.1535
$ws.Range("M134").Value = -4652.1921
$ws.Range("N134").Value = -13798.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 165012.75
$ws.Range("J2").Value = 13.285714
$ws.Range("L2").Value = 79.71428400000001
$ws.Range("N2").Value = -305.714284

$ws.Range("H12").Value = 53.115383
$ws.Range("I12").Value = 95.666664
$ws.Range("K12").Value = 286.999992
$ws.Range("M12").Value = -113.999992

$ws.Range("H23").Value = 90
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 270
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -740

$ws.Range("H25").Value = 1320.4
$ws.Range("J25").Value = 1320.4
$ws.Range("L25").Value = 3961.2
$ws.Range("N25").Value = -4299.200000000001

$ws.Range("H30").Value = 1320.4
$ws.Range("J30").Value = 1320.4
$ws.Range("L30").Value = 3961.2
$ws.Range("N30").Value = -4165.200000000001

$ws.Range("H132").Value = 1502.9048
$ws.Range("J132").Value = 1533.05
$ws.Range("L132").Value = 13797.45
$ws.Range("N132").Value = -18857.45

$ws.Range("H141").Value = 1485.9445
$ws.Range("I141").Value = 821
$ws.Range("J141").Value = 1909.091
$ws.Range("K141").Value = 2463
$ws.Range("L141").Value = 5727.272999999999
$ws.Range("M141").Value = 2717
$ws.Range("N141").Value = -16087.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2565
$ws.Range("I7").Value = 1700
$ws.Range("J7").Value = 3141.6667
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 3141.6667
$ws.Range("M7").Value = -1588
$ws.Range("N7").Value = -3365.6667

$ws.Range("H40").Value = 2200.8
$ws.Range("I40").Value = 1452
$ws.Range("J40").Value = 2700
$ws.Range("K40").Value = 1452
$ws.Range("L40").Value = 2700
$ws.Range("M40").Value = -1316
$ws.Range("N40").Value = -2972

$ws.Range("H61").Value = 3861.8
$ws.Range("I61").Value = 4652
$ws.Range("J61").Value = 3335
$ws.Range("K61").Value = 4652
$ws.Range("L61").Value = 3335
$ws.Range("M61").Value = -4450
$ws.Range("N61").Value = -3739

$ws.Range("H68").Value = 13535395
$ws.Range("I68").Value = 28195482
$ws.Range("J68").Value = 3008
$ws.Range("K68").Value = 28195482
$ws.Range("L68").Value = 3008
$ws.Range("M68").Value = -28194733
$ws.Range("N68").Value = -4506

$ws.Range("H71").Value = 13535395
$ws.Range("I71").Value = 28195482
$ws.Range("J71").Value = 3008
$ws.Range("K71").Value = 140977410
$ws.Range("L71").Value = 15040
$ws.Range("M71").Value = -140973666
$ws.Range("N71").Value = -22528

$ws.Range("H82").Value = 6495571
$ws.Range("I82").Value = 11364574
$ws.Range("K82").Value = 11364574
$ws.Range("M82").Value = -11364213

$ws.Range("H85").Value = 6495571
$ws.Range("I85").Value = 11364574
$ws.Range("K85").Value = 11364574
$ws.Range("M85").Value = -11363326

$ws.Range("H113").Value = 3861.8
$ws.Range("I113").Value = 4652
$ws.Range("J113").Value = 3335
$ws.Range("K113").Value = 4652
$ws.Range("L113").Value = 3335
$ws.Range("M113").Value = -2482
$ws.Range("N113").Value = -7675

$ws.Range("H122").Value = 4195.8887
$ws.Range("I122").Value = 4599.5264
$ws.Range("J122").Value = 3237.25
$ws.Range("K122").Value = 13798.5792
$ws.Range("L122").Value = 9711.75
$ws.Range("M122").Value = -11348.5792
$ws.Range("N122").Value = -14611.75

$ws.Range("H126").Value = 2565
$ws.Range("I126").Value = 1700
$ws.Range("J126").Value = 3141.6667
$ws.Range("K126").Value = 5100
$ws.Range("L126").Value = 9425.000100000001
$ws.Range("M126").Value = -2630
$ws.Range("N126").Value = -14365.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 76925830
$ws.Range("I81").Value = 166669460
$ws.Range("K81").Value = 333338920
$ws.Range("M81").Value = -333337859

$ws.Range("H84").Value = 76925830
$ws.Range("I84").Value = 166669460
$ws.Range("K84").Value = 1666694600
$ws.Range("M84").Value = -1666689296

$ws.Range("H107").Value = 1391.0435
$ws.Range("I107").Value = 1171.5
$ws.Range("J107").Value = 1732.5555
$ws.Range("K107").Value = 3514.5
$ws.Range("L107").Value = 5197.666499999999
$ws.Range("M107").Value = -1594.5
$ws.Range("N107").Value = -9037.666499999999

$ws.Range("H113").Value = 100000660
